# Applies the "Fit_distribs y probs" edit to the Data sheet:
#  - Adds a new column C ("Calama 2") mirroring column B's ID/label rows,
#    and for the data rows (5-58) copies the numeric value from column B
#    when it is a real observation (>= 1), or records #N/A when the value
#    produced no match (i.e. B was 0 or < 1), matching the probs lookup
#    results.
#  - Formats: C1/C4 reuse the header fill style of B1/B4, C2 reuses the
#    highlighted style of B2, and the C5:C58 data block is center aligned.
#  - Updates the sheet selection and page setup to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header row (label row referencing B4 header) ---
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1").Value = "Calama 2"

# --- Row 2: "Datos" row, highlighted count cell ---
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = 54

# --- Row 4: second header row (actual column title) ---
$ws.Range("B4").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$ws.Range("C4").Value = "Calama 2"

# --- Rows 5-58: data block, center aligned ---
$ws.Range("C5:C58").HorizontalAlignment = -4108   # xlCenter

$dataMap = @{
    5  = 6.7;    6  = "#N/A"; 7  = 2.7;    8  = 1.4;    9  = "#N/A";
    10 = "#N/A"; 11 = "#N/A"; 12 = 9.5;    13 = "#N/A"; 14 = "#N/A";
    15 = 7;      16 = 6;      17 = 1.5;    18 = "#N/A"; 19 = 4;
    20 = "#N/A"; 21 = "#N/A"; 22 = "#N/A"; 23 = 13.5;   24 = 5.2;
    25 = "#N/A"; 26 = 1.5;    27 = 3;      28 = 3.1;    29 = 2;
    30 = "#N/A"; 31 = "#N/A"; 32 = 7.5;    33 = "#N/A"; 34 = "#N/A";
    35 = "#N/A"; 36 = 3.5;    37 = "#N/A"; 38 = "#N/A"; 39 = "#N/A";
    40 = 1.5;    41 = 7.7;    42 = 3.7;    43 = 7.5;    44 = "#N/A";
    45 = 3.5;    46 = "#N/A"; 47 = "#N/A"; 48 = "#N/A"; 49 = "#N/A";
    50 = "#N/A"; 51 = 4;      52 = 4;      53 = 2;      54 = 4.4;
    55 = 7.4;    56 = "#N/A"; 57 = 1.3;    58 = 3.5
}

foreach ($row in 5..58) {
    $ws.Range("C$row").Value = $dataMap[$row]
}

# --- Sheet selection / page setup to match saved state ---
[void]$ws.Range("D3").Select()
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Edit applied"
